$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(3.3693521517586755, 1.6592021395484187, -1.546351642382673, 2.614870506703698)"
$ws.Range("C2").Value = "JSU(-1.2578714585130992, 1.193165094177969, 2.960283936233231, 5.153206108507694)"
$ws.Range("D2").Value = "NIG(0.7160223505976502, 0.530703991212947, 1.1795464246861718, 2.810590625915288)"
$ws.Range("E2").Value = "JSU(-1.3560890004458988, 1.3008065720638309, 1.695757897553844, 5.551220762467031)"
